# Sprint #1 Burndown Chart - update Actual Work burndown for the final
# working day (idx 14 / row 16) and set the sheet's print/page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (day 14) "Actual Work" (column C) drops from 67 - 5.5 = 61.5
# down to 67 - 10.5 = 56.5 -- an extra 5 points of work burned down.
$ws.Range("C16").Formula = "=C15-10.5"

# Documentation: set the worksheet up for A4 portrait printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
